$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")

Set-TextValue $ws1.Range("C2") "$2,057.00"
Set-TextValue $ws1.Range("C3") "$951.97"
Set-TextValue $ws1.Range("C5") "$1,750.00 USD"
$ws1.Range("A6").Value = "XL™ HALF RACK"
Set-TextValue $ws1.Range("C6") "$2,399.00"
$ws1.Range("F6").Value = "https://www.sorinex.com/products/xl-half-rack?Attachment+Color=Black+Texture&Upgrades=None"

$ws2 = $wb.Worksheets.Item("Squat Stands")
Set-TextValue $ws2.Range("C2") "$1,488.00"
